# Daily attendance processing - normalize "Recorded By" (column G) entries so
# that any leading "System" / "system" author tag is moved to the end of the
# comma-separated list (swap the first and last entries), with the moved
# "System" token re-capitalized to "System".
#
# Example: "System, dnasr281@gmail.com"          -> "dnasr281@gmail.com, System"
#          "System, backup@backdoor.com, system" -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G = "Recorded By" (header is in row 1).
$col = 7

# Find the last used row in column A (xlUp from the bottom of the sheet).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ', '

    if ($parts.Count -gt 1 -and $parts[0].ToLower() -eq 'system') {
        $lastIdx = $parts.Count - 1
        $firstOriginal = $parts[0]
        $parts[0] = $parts[$lastIdx]
        $parts[$lastIdx] = 'System'
        $newVal = $parts -join ', '
        $cell.Value = $newVal
    }
}
